# Updated symbol list on Tue Jan 31 03:38:57 UTC 2023 with GitHub Actions
# Refresh latest crypto price/volume figures on the "cryptos" sheet.
# Values are stored as plain text (matching the sheet's existing inline-string
# convention), so Text number format is forced before assigning the new
# value and the default "Normal" style is restored afterward so no new
# formatting/style gets introduced on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = "Normal"
}

Set-TextValue "E2" "-1.64%"
Set-TextValue "D3" "38.11"
Set-TextValue "E3" "-3.38%"
Set-TextValue "D4" "5.067"
Set-TextValue "E4" "-1.28%"
Set-TextValue "D5" "0.07762"
Set-TextValue "E5" "-4.97%"
Set-TextValue "D6" "4.357"
Set-TextValue "E6" "-0.26%"
Set-TextValue "D7" "1.901"
Set-TextValue "E7" "-4.19%"
Set-TextValue "D8" "8.194"
Set-TextValue "E8" "-1.67%"
Set-TextValue "D9" "0.9221"
Set-TextValue "E9" "-1.57%"
Set-TextValue "E10" "-4.80%"
Set-TextValue "D11" "0.1879"
Set-TextValue "E11" "-4.65%"
Set-TextValue "D12" "0.08764"
Set-TextValue "E12" "-2.96%"
Set-TextValue "D13" "0.03421"
Set-TextValue "E13" "-2.19%"
Set-TextValue "D14" "0.09709"
Set-TextValue "E14" "-0.38%"
Set-TextValue "D15" "0.001377"
Set-TextValue "E15" "-2.63%"
Set-TextValue "D16" "0.005730"
Set-TextValue "E16" "-4.72%"
Set-TextValue "D17" "3.563"
Set-TextValue "E17" "-2.15%"
Set-TextValue "E18" "-6.50%"
Set-TextValue "E19" "-2.28%"
Set-TextValue "D20" "5.026"
Set-TextValue "E20" "1.40%"
Set-TextValue "E21" "-3.76%"
Set-TextValue "D22" "0.2620"
Set-TextValue "E22" "1.53%"
Set-TextValue "E23" "5,592.34%"
Set-TextValue "D24" "0.04392"
Set-TextValue "E24" "0.36%"
Set-TextValue "E25" "-2.59%"
Set-TextValue "D26" "0.004255"
Set-TextValue "E26" "-10.68%"
Set-TextValue "E27" "-65.30%"
Set-TextValue "D39" "0.02133"
Set-TextValue "E39" "-3.71%"
Set-TextValue "D40" "0.05008"
Set-TextValue "E40" "-3.39%"
Set-TextValue "D41" "0.007808"
Set-TextValue "E41" "0.44%"
Set-TextValue "D42" "0.01005"
Set-TextValue "E42" "-2.08%"
Set-TextValue "D43" "0.1343"
Set-TextValue "E43" "-4.00%"
Set-TextValue "D44" "0.002061"
Set-TextValue "E44" "-1.93%"
Set-TextValue "D45" "0.008799"
Set-TextValue "E45" "-5.15%"
Set-TextValue "E46" "-6.71%"
Set-TextValue "E47" "-0.02%"
Set-TextValue "D48" "0.003246"
Set-TextValue "E48" "12.58%"
Set-TextValue "E49" "-0.13%"
Set-TextValue "E50" "-0.02%"
Set-TextValue "E51" "-0.02%"
